$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "(318294931, Shalev  Afanasenko: -2,-9)"
$ws.Range("B1").Value = "(305487936, Avihai  Kipnis: 0,8)"
$ws.Range("C1").Value = "(313227928, Aviv  Levi: 6,5)"
$ws.Range("D1").Value = "(205807308, Sariel  Basis: 2,4)"
$ws.Range("E1").Value = "(315891549, Raz  Halaby: 2,-8)"
$ws.Range("F1").Value = "(315060103, Dan  Mshelh: -10,-5)"
$ws.Range("G1").Value = "(313925141, Elad   Amer: -1,6)"

$ws.Range("A3").Value = "cost: 475.6196525797319"
$ws.Range("A4").Value = "time: 64.37423608281884"
